{"js": "// Update the \"three-digit number multiplied by one-digit number\" answer\n// table: each existing \"A\u00d7B=C\" cell text is replaced by a new equation,\n// one-for-one, in document order. Matching is done via exact text search\n// so the run's existing formatting (font/size) is preserved untouched.\nconst replacements = [\n  [\"692\u00d76=4152\", \"692\u00d72=1384\"],\n  [\"413\u00d78=3304\", \"305\u00d72=610\"],\n  [\"369\u00d75=1845\", \"128\u00d76=768\"],\n  [\"765\u00d76=4590\", \"762\u00d76=4572\"],\n  [\"427\u00d76=2562\", \"374\u00d78=2992\"],\n  [\"153\u00d72=306\", \"856\u00d72=1712\"],\n  [\"297\u00d79=2673\", \"502\u00d72=1004\"],\n  [\"728\u00d74=2912\", \"753\u00d73=2259\"],\n  [\"712\u00d77=4984\", \"237\u00d77=1659\"],\n  [\"885\u00d79=7965\", \"414\u00d74=1656\"],\n  [\"446\u00d73=1338\", \"722\u00d76=4332\"],\n  [\"652\u00d76=3912\", \"705\u00d78=5640\"],\n  [\"803\u00d77=5621\", \"473\u00d76=2838\"],\n  [\"151\u00d79=1359\", \"702\u00d76=4212\"],\n  [\"731\u00d78=5848\", \"308\u00d75=1540\"],\n  [\"209\u00d79=1881\", \"754\u00d78=6032\"],\n  [\"438\u00d78=3504\", \"888\u00d72=1776\"],\n  [\"571\u00d79=5139\", \"739\u00d75=3695\"],\n  [\"320\u00d72=640\", \"772\u00d72=1544\"],\n  [\"898\u00d75=4490\", \"863\u00d73=2589\"],\n  [\"578\u00d72=1156\", \"388\u00d72=776\"],\n  [\"491\u00d77=3437\", \"613\u00d78=4904\"],\n  [\"234\u00d72=468\", \"616\u00d78=4928\"],\n  [\"646\u00d79=5814\", \"342\u00d76=2052\"],\n  [\"818\u00d75=4090\", \"456\u00d73=1368\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  // Replace the first (and expected-only) occurrence, preserving the\n  // run's existing formatting.\n  found.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update the \"three-digit number multiplied by one-digit number\" answer\n# table: each existing \"A\u00d7B=C\" cell text is replaced by a new equation,\n# one-for-one, in document order. Find/Replace is scoped to an exact,\n# case-sensitive match of the whole old equation string so the\n# surrounding run formatting (font/size) is left untouched.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{ Old = \"692\u00d76=4152\"; New = \"692\u00d72=1384\" },\n  @{ Old = \"413\u00d78=3304\"; New = \"305\u00d72=610\" },\n  @{ Old = \"369\u00d75=1845\"; New = \"128\u00d76=768\" },\n  @{ Old = \"765\u00d76=4590\"; New = \"762\u00d76=4572\" },\n  @{ Old = \"427\u00d76=2562\"; New = \"374\u00d78=2992\" },\n  @{ Old = \"153\u00d72=306\"; New = \"856\u00d72=1712\" },\n  @{ Old = \"297\u00d79=2673\"; New = \"502\u00d72=1004\" },\n  @{ Old = \"728\u00d74=2912\"; New = \"753\u00d73=2259\" },\n  @{ Old = \"712\u00d77=4984\"; New = \"237\u00d77=1659\" },\n  @{ Old = \"885\u00d79=7965\"; New = \"414\u00d74=1656\" },\n  @{ Old = \"446\u00d73=1338\"; New = \"722\u00d76=4332\" },\n  @{ Old = \"652\u00d76=3912\"; New = \"705\u00d78=5640\" },\n  @{ Old = \"803\u00d77=5621\"; New = \"473\u00d76=2838\" },\n  @{ Old = \"151\u00d79=1359\"; New = \"702\u00d76=4212\" },\n  @{ Old = \"731\u00d78=5848\"; New = \"308\u00d75=1540\" },\n  @{ Old = \"209\u00d79=1881\"; New = \"754\u00d78=6032\" },\n  @{ Old = \"438\u00d78=3504\"; New = \"888\u00d72=1776\" },\n  @{ Old = \"571\u00d79=5139\"; New = \"739\u00d75=3695\" },\n  @{ Old = \"320\u00d72=640\"; New = \"772\u00d72=1544\" },\n  @{ Old = \"898\u00d75=4490\"; New = \"863\u00d73=2589\" },\n  @{ Old = \"578\u00d72=1156\"; New = \"388\u00d72=776\" },\n  @{ Old = \"491\u00d77=3437\"; New = \"613\u00d78=4904\" },\n  @{ Old = \"234\u00d72=468\"; New = \"616\u00d78=4928\" },\n  @{ Old = \"646\u00d79=5814\"; New = \"342\u00d76=2052\" },\n  @{ Old = \"818\u00d75=4090\"; New = \"456\u00d73=1368\" }\n)\n\nforeach ($p in $pairs) {\n  $r = $d.Content\n  $r.Find.ClearFormatting()\n  $r.Find.Replacement.ClearFormatting()\n  $r.Find.Text = $p.Old\n  $r.Find.Replacement.Text = $p.New\n  $r.Find.Forward = $true\n  $r.Find.Wrap = 1\n  $r.Find.Format = $false\n  $r.Find.MatchCase = $true\n  $r.Find.MatchWholeWord = $false\n  $r.Find.MatchWildcards = $false\n  $r.Find.MatchSoundsLike = $false\n  $r.Find.MatchAllWordForms = $false\n  $found = $r.Find.Execute($r.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $r.Find.Replacement.Text, 2)\n  if (-not $found) {\n    Write-Output (\"NOT FOUND: \" + $p.Old)\n  }\n}\n"}
